$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the year values in column A (payment month/year reference values)
$ws.Range("A2").Value = 2020
$ws.Range("A3").Value = 2019
$ws.Range("A4").Value = 2018

# Set explicit column widths for columns B and C
$ws.Columns.Item(2).ColumnWidth = 10.33
$ws.Columns.Item(3).ColumnWidth = 11.14

# Move the active selection to D4
$ws.Range("D4").Select()
